$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "PopularCarModels" worksheet after the existing last sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PopularCarModels"

# Header + popular car model names
$values = @(
    "Popular Car Model",
    "Maruti 800",
    "Maruti Swift Dzire",
    "Maruti Swift",
    "Hyundai I10",
    "Hyundai Santro Xing",
    "Honda City",
    "Toyota Innova",
    "Toyota Fortuner",
    "Mahindra XUV500"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Style the header cell like the header row of BikeDetails (white text on the green fill)
$header = $ws.Range("A1")
$header.Font.ColorIndex = 2
$header.Interior.ColorIndex = 10
$header.Interior.Pattern = 1

# Size column A to fit the content, matching the BikeDetails header column sizing convention
$ws.Columns.Item(1).ColumnWidth = 16.34375

# Restore BikeDetails as the active/selected sheet
$ws1.Activate()
